$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value changes (leading apostrophe preserves the existing
#     quotePrefix cell style instead of Excel picking a different xf) ---

# Row 2 (order matches the new shared-string append order)
$ws.Range("V2").Value = "'EMEAAD\abouhadjer"
$ws.Range("Q2").Value = "'s"
$ws.Range("U2").Value = "'EMEAAD\alepicard"
$ws.Range("O2").Value = "'HW - Network Security"
$ws.Range("R2").Value = "'MAINT FIXED COST"

# Row 3
$ws.Range("U3").Value = "'"

# --- Column width changes ---
$ws.Range("R1").ColumnWidth = 17.7109375
$ws.Range("U1").ColumnWidth = 19.85546875

# --- Sheet view changes ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$ws.Range("R2").Select()
